$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped from 45203 to 45204
# (one day later) for every data row (rows 2 through 490).
$ws.Range("C2:C490").Value = 45204
